# Natmi following Dr Hou advice
# Rebuild the Sending-cluster x Target-cluster cross table for L1cam-Ephb2
# (ECs/FAPs/sCs x ECs/FAPs/sCs) with updated NATMI statistics, including the
# previously-missing self-pairs (ECs->ECs, FAPs->FAPs, sCs->sCs) as rows 8-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "L1cam"
$ws.Cells.Item(2,3).Value = "Ephb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 22.59487733333333
$ws.Cells.Item(2,8).Value = 67.784632
$ws.Cells.Item(2,9).Value = 0.7395019553569895
$ws.Cells.Item(2,10).Value = 0.7395019553569895
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.005966
$ws.Cells.Item(2,14).Value = 0.017898
$ws.Cells.Item(2,15).Value = 0.00125947234315407
$ws.Cells.Item(2,16).Value = 0.00125947234315407
$ws.Cells.Item(2,17).Value = 0.1348010381706667
$ws.Cells.Item(2,18).Value = 1.213209343536
$ws.Cells.Item(2,19).Value = 0.0009313822604804839
$ws.Cells.Item(2,20).Value = 0.0009313822604804839

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "L1cam"
$ws.Cells.Item(3,3).Value = "Ephb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 22.59487733333333
$ws.Cells.Item(3,8).Value = 67.784632
$ws.Cells.Item(3,9).Value = 0.7395019553569895
$ws.Cells.Item(3,10).Value = 0.7395019553569895
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.821776
$ws.Cells.Item(3,14).Value = 11.465328
$ws.Cells.Item(3,15).Value = 0.8068087787009701
$ws.Cells.Item(3,16).Value = 0.8068087787009701
$ws.Cells.Item(3,17).Value = 86.35255991547733
$ws.Cells.Item(3,18).Value = 777.173039239296
$ws.Cells.Item(3,19).Value = 0.596636669448552
$ws.Cells.Item(3,20).Value = 0.596636669448552

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "L1cam"
$ws.Cells.Item(4,3).Value = "Ephb2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 22.59487733333333
$ws.Cells.Item(4,8).Value = 67.784632
$ws.Cells.Item(4,9).Value = 0.7395019553569895
$ws.Cells.Item(4,10).Value = 0.7395019553569895
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.9091623333333333
$ws.Cells.Item(4,14).Value = 2.727487
$ws.Cells.Item(4,15).Value = 0.1919317489558758
$ws.Cells.Item(4,16).Value = 0.1919317489558758
$ws.Cells.Item(4,17).Value = 20.54241139775378
$ws.Cells.Item(4,18).Value = 184.881702579784
$ws.Cells.Item(4,19).Value = 0.1419339036479569
$ws.Cells.Item(4,20).Value = 0.1419339036479569

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "L1cam"
$ws.Cells.Item(5,3).Value = "Ephb2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.3045986666666667
$ws.Cells.Item(5,8).Value = 0.913796
$ws.Cells.Item(5,9).Value = 0.00996913177602551
$ws.Cells.Item(5,10).Value = 0.00996913177602551
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.005966
$ws.Cells.Item(5,14).Value = 0.017898
$ws.Cells.Item(5,15).Value = 0.00125947234315407
$ws.Cells.Item(5,16).Value = 0.00125947234315407
$ws.Cells.Item(5,17).Value = 0.001817235645333333
$ws.Cells.Item(5,18).Value = 0.016355120808
$ws.Cells.Item(5,19).Value = 0.00001255584575716254
$ws.Cells.Item(5,20).Value = 0.00001255584575716254

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "L1cam"
$ws.Cells.Item(6,3).Value = "Ephb2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.3045986666666667
$ws.Cells.Item(6,8).Value = 0.913796
$ws.Cells.Item(6,9).Value = 0.00996913177602551
$ws.Cells.Item(6,10).Value = 0.00996913177602551
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.821776
$ws.Cells.Item(6,14).Value = 11.465328
$ws.Cells.Item(6,15).Value = 0.8068087787009701
$ws.Cells.Item(6,16).Value = 0.8068087787009701
$ws.Cells.Item(6,17).Value = 1.164107873898667
$ws.Cells.Item(6,18).Value = 10.476970865088
$ws.Cells.Item(6,19).Value = 0.008043183032924174
$ws.Cells.Item(6,20).Value = 0.008043183032924174

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "L1cam"
$ws.Cells.Item(7,3).Value = "Ephb2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.3045986666666667
$ws.Cells.Item(7,8).Value = 0.913796
$ws.Cells.Item(7,9).Value = 0.00996913177602551
$ws.Cells.Item(7,10).Value = 0.00996913177602551
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.9091623333333333
$ws.Cells.Item(7,14).Value = 2.727487
$ws.Cells.Item(7,15).Value = 0.1919317489558758
$ws.Cells.Item(7,16).Value = 0.1919317489558758
$ws.Cells.Item(7,17).Value = 0.2769296345168889
$ws.Cells.Item(7,18).Value = 2.492366710652
$ws.Cells.Item(7,19).Value = 0.001913392897344172
$ws.Cells.Item(7,20).Value = 0.001913392897344172

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "L1cam"
$ws.Cells.Item(8,3).Value = "Ephb2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 7.654706
$ws.Cells.Item(8,8).Value = 22.964118
$ws.Cells.Item(8,9).Value = 0.2505289128669849
$ws.Cells.Item(8,10).Value = 0.2505289128669849
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.005966
$ws.Cells.Item(8,14).Value = 0.017898
$ws.Cells.Item(8,15).Value = 0.00125947234315407
$ws.Cells.Item(8,16).Value = 0.00125947234315407
$ws.Cells.Item(8,17).Value = 0.045667975996
$ws.Cells.Item(8,18).Value = 0.411011783964
$ws.Cells.Item(8,19).Value = 0.0003155342369164233
$ws.Cells.Item(8,20).Value = 0.0003155342369164233

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "L1cam"
$ws.Cells.Item(9,3).Value = "Ephb2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 7.654706
$ws.Cells.Item(9,8).Value = 22.964118
$ws.Cells.Item(9,9).Value = 0.2505289128669849
$ws.Cells.Item(9,10).Value = 0.2505289128669849
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.821776
$ws.Cells.Item(9,14).Value = 11.465328
$ws.Cells.Item(9,15).Value = 0.8068087787009701
$ws.Cells.Item(9,16).Value = 0.8068087787009701
$ws.Cells.Item(9,17).Value = 29.254571677856
$ws.Cells.Item(9,18).Value = 263.291145100704
$ws.Cells.Item(9,19).Value = 0.2021289262194939
$ws.Cells.Item(9,20).Value = 0.2021289262194939

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "L1cam"
$ws.Cells.Item(10,3).Value = "Ephb2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.654706
$ws.Cells.Item(10,8).Value = 22.964118
$ws.Cells.Item(10,9).Value = 0.2505289128669849
$ws.Cells.Item(10,10).Value = 0.2505289128669849
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.9091623333333333
$ws.Cells.Item(10,14).Value = 2.727487
$ws.Cells.Item(10,15).Value = 0.1919317489558758
$ws.Cells.Item(10,16).Value = 0.1919317489558758
$ws.Cells.Item(10,17).Value = 6.959370367940666
$ws.Cells.Item(10,18).Value = 62.63433331146599
$ws.Cells.Item(10,19).Value = 0.04808445241057462
$ws.Cells.Item(10,20).Value = 0.04808445241057462
